$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 81; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value()
    $oa = $current.ToOADate()
    if ($oa -eq 45243) {
        $cell.Value = 45244
    }
}
